# Scheduled-runner refresh of cached market/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job Leve sheets. Values below were
# recomputed from a fresh Universalis price pull; cells that no longer carry a
# NQ/HQ split are cleared (value "") rather than left stale, and newly-split
# rows gain their profit cell back.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 490
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = ""
$ws.Range("H137").Value = 3031914.5
$ws.Range("I137").Value = 5883391
$ws.Range("J137").Value = 2220.6875
$ws.Range("K137").Value = 17650173
$ws.Range("L137").Value = 6662.0625
$ws.Range("M137").Value = -17647623
$ws.Range("N137").Value = -11762.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18785.828
$ws.Range("I32").Value = 19321.018
$ws.Range("J32").Value = 15915.272
$ws.Range("K32").Value = 19321.018
$ws.Range("L32").Value = 15915.272
$ws.Range("M32").Value = -19034.018
$ws.Range("N32").Value = -16489.272
$ws.Range("H61").Value = 222456400
$ws.Range("I61").Value = 166834600
$ws.Range("J61").Value = 333700000
$ws.Range("K61").Value = 166834600
$ws.Range("L61").Value = 333700000
$ws.Range("M61").Value = -166834388
$ws.Range("N61").Value = -333700424
$ws.Range("H74").Value = 12601017
$ws.Range("I74").Value = 15688528
$ws.Range("J74").Value = 250971
$ws.Range("K74").Value = 15688528
$ws.Range("L74").Value = 250971
$ws.Range("M74").Value = -15687654
$ws.Range("N74").Value = -252719
$ws.Range("H77").Value = 12601017
$ws.Range("I77").Value = 15688528
$ws.Range("J77").Value = 250971
$ws.Range("K77").Value = 78442640
$ws.Range("L77").Value = 1254855
$ws.Range("M77").Value = -78438272
$ws.Range("N77").Value = -1263591
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = ""
$ws.Range("H109").Value = 25000
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27774
$ws.Range("H132").Value = 51975.6
$ws.Range("I132").Value = 39560.23
$ws.Range("J132").Value = 75032.71000000001
$ws.Range("K132").Value = 118680.69
$ws.Range("L132").Value = 225098.13
$ws.Range("M132").Value = -116150.69
$ws.Range("N132").Value = -230158.13
$ws.Range("H136").Value = 222456400
$ws.Range("I136").Value = 166834600
$ws.Range("J136").Value = 333700000
$ws.Range("K136").Value = 500503800
$ws.Range("L136").Value = 1001100000
$ws.Range("M136").Value = -500501250
$ws.Range("N136").Value = -1001105100
$ws.Range("H141").Value = 45000
$ws.Range("J141").Value = 45000
$ws.Range("L141").Value = 45000
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2644.3333
$ws.Range("I99").Value = 2300
$ws.Range("J99").Value = 3333
$ws.Range("K99").Value = 2300
$ws.Range("L99").Value = 3333
$ws.Range("M99").Value = -802
$ws.Range("N99").Value = -6329
$ws.Range("H107").Value = 3277.6155
$ws.Range("I107").Value = 3892.5
$ws.Range("J107").Value = 3004.3333
$ws.Range("K107").Value = 3892.5
$ws.Range("L107").Value = 3004.3333
$ws.Range("M107").Value = -1972.5
$ws.Range("N107").Value = -6844.3333
$ws.Range("H134").Value = 2760.3809
$ws.Range("I134").Value = 3016.125
$ws.Range("J134").Value = 1942
$ws.Range("K134").Value = 9048.375
$ws.Range("L134").Value = 5826
$ws.Range("M134").Value = -6513.375
$ws.Range("N134").Value = -10896
$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1511.0834
$ws.Range("I31").Value = 1193.909
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1193.909
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -898.9090000000001
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 1511.0834
$ws.Range("I34").Value = 1193.909
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1193.909
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -991.9090000000001
$ws.Range("N34").Value = -5404
$ws.Range("H58").Value = 45456130
$ws.Range("I58").Value = 35715264
$ws.Range("J58").Value = 62502644
$ws.Range("K58").Value = 35715264
$ws.Range("L58").Value = 62502644
$ws.Range("M58").Value = -35715061
$ws.Range("N58").Value = -62503050
$ws.Range("H106").Value = 11491.6
$ws.Range("J106").Value = 11491.6
$ws.Range("L106").Value = 11491.6
$ws.Range("N106").Value = -14015.6
$ws.Range("H132").Value = 32027.152
$ws.Range("I132").Value = 1544.2858
$ws.Range("J132").Value = 85372.164
$ws.Range("K132").Value = 4632.857400000001
$ws.Range("L132").Value = 256116.492
$ws.Range("M132").Value = -2102.857400000001
$ws.Range("N132").Value = -261176.492
$ws.Range("H134").Value = 49099
$ws.Range("I134").Value = 2339.158
$ws.Range("J134").Value = 271208.25
$ws.Range("K134").Value = 7017.474
$ws.Range("L134").Value = 813624.75
$ws.Range("M134").Value = -4482.474
$ws.Range("N134").Value = -818694.75
$ws.Range("H136").Value = 45456130
$ws.Range("I136").Value = 35715264
$ws.Range("J136").Value = 62502644
$ws.Range("K136").Value = 107145792
$ws.Range("L136").Value = 187507932
$ws.Range("M136").Value = -107143242
$ws.Range("N136").Value = -187513032

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1848.2727
$ws.Range("I97").Value = 1565
$ws.Range("J97").Value = 2603.6667
$ws.Range("K97").Value = 1565
$ws.Range("L97").Value = 2603.6667
$ws.Range("M97").Value = -1069
$ws.Range("N97").Value = -3595.6667
$ws.Range("H132").Value = 61738.56
$ws.Range("I132").Value = 47736.953
$ws.Range("J132").Value = 87408.164
$ws.Range("K132").Value = 143210.859
$ws.Range("L132").Value = 262224.492
$ws.Range("M132").Value = -140680.859
$ws.Range("N132").Value = -267284.492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 92172.73
$ws.Range("I100").Value = 92172.73
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 92172.73
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -91631.73
$ws.Range("N100").Value = ""
$ws.Range("H132").Value = 128310
$ws.Range("I132").Value = 2166.6667
$ws.Range("J132").Value = 203996
$ws.Range("K132").Value = 6500.000100000001
$ws.Range("L132").Value = 611988
$ws.Range("M132").Value = -3970.000100000001
$ws.Range("N132").Value = -617048
$ws.Range("H136").Value = 118759.35
$ws.Range("I136").Value = 167566.5
$ws.Range("J136").Value = 92137.27
$ws.Range("K136").Value = 502699.5
$ws.Range("L136").Value = 276411.81
$ws.Range("M136").Value = -500149.5
$ws.Range("N136").Value = -281511.81

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 13706.857
$ws.Range("J54").Value = 13706.857
$ws.Range("L54").Value = 13706.857
$ws.Range("N54").Value = -14746.857
$ws.Range("H81").Value = 2535.0588
$ws.Range("I81").Value = 1645
$ws.Range("J81").Value = 2653.7334
$ws.Range("K81").Value = 3290
$ws.Range("L81").Value = 5307.4668
$ws.Range("M81").Value = -2229
$ws.Range("N81").Value = -7429.4668
$ws.Range("H84").Value = 2535.0588
$ws.Range("I84").Value = 1645
$ws.Range("J84").Value = 2653.7334
$ws.Range("K84").Value = 16450
$ws.Range("L84").Value = 26537.334
$ws.Range("M84").Value = -11146
$ws.Range("N84").Value = -37145.334
$ws.Range("H132").Value = 107060.31
$ws.Range("I132").Value = 84265.664
$ws.Range("J132").Value = 146136.86
$ws.Range("K132").Value = 252796.992
$ws.Range("L132").Value = 438410.58
$ws.Range("M132").Value = -250266.992
$ws.Range("N132").Value = -443470.58
$ws.Range("H136").Value = 47577.836
$ws.Range("I136").Value = 26626.744
$ws.Range("J136").Value = 251851
$ws.Range("K136").Value = 79880.23199999999
$ws.Range("L136").Value = 755553
$ws.Range("M136").Value = -77330.23199999999
$ws.Range("N136").Value = -760653
